# Fruta / hortaliza, semanal
# Insert 3 new weekly price rows for "Femacal de La Calera - Frutilla" (row 237)
# which shifts the existing rows 237-344 down to 240-347.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new blank rows before row 237 (existing rows 237.. shift down by 3)
$ws.Rows.Item(237).Resize(3).Insert()

# Data that is constant across every row of this table
$mercadoId = 3
$mercado   = "Femacal de La Calera"
$region    = "Coquimbo"
$codreg    = 5
$tipo      = "Fruta"
$productoId = 100101
$producto   = "Berries"
$categoriaId = 100112025
$categoria   = "Frutilla"
$variedad    = "Sin especificar"
$unidad      = "`$/bandeja 7 kilos"
$origen      = "Provincia de Melipilla"
$kgUnidad    = 7

# New rows (237, 238, 239): Fecha, Calidad, Volumen, PrecioMin, PrecioMax, PrecioProm, PrecioKg
$newRows = @(
    @{ Row = 237; Fecha = 44845; Calidad = "Especial"; Volumen = 50; PMin = 15000; PMax = 15000; PProm = 15000; PKg = 2143 },
    @{ Row = 238; Fecha = 44845; Calidad = "Primera";  Volumen = 57; PMin = 12000; PMax = 12000; PProm = 12000; PKg = 1714 },
    @{ Row = 239; Fecha = 44845; Calidad = "Segunda";  Volumen = 45; PMin = 9000;  PMax = 9000;  PProm = 9000;  PKg = 1286 }
)

foreach ($r in $newRows) {
    $row = $r.Row

    $ws.Cells.Item($row, 1).Value = $mercadoId
    $ws.Cells.Item($row, 2).Value = $mercado
    $ws.Cells.Item($row, 3).Value = $region
    $ws.Cells.Item($row, 4).Value = $r.Fecha
    $ws.Cells.Item($row, 5).Value = $codreg
    $ws.Cells.Item($row, 6).Value = $tipo
    $ws.Cells.Item($row, 7).Value = $productoId
    $ws.Cells.Item($row, 8).Value = $producto
    $ws.Cells.Item($row, 9).Value = $categoriaId
    $ws.Cells.Item($row, 10).Value = $categoria
    $ws.Cells.Item($row, 11).Value = $variedad
    $ws.Cells.Item($row, 12).Value = $r.Calidad
    $ws.Cells.Item($row, 13).Value = $r.Volumen
    $ws.Cells.Item($row, 14).Value = $r.PMin
    $ws.Cells.Item($row, 15).Value = $r.PMax
    $ws.Cells.Item($row, 16).Value = $r.PProm
    $ws.Cells.Item($row, 17).Value = $unidad
    $ws.Cells.Item($row, 18).Value = $origen
    $ws.Cells.Item($row, 19).Value = $r.PKg
    $ws.Cells.Item($row, 20).Value = $kgUnidad
}
